$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 30   Number  34"
$ws.Range("C9").Value = "Report Covering the Week  8/21/2023  Through  8/27/2023"

# --- Cells that convert from numeric to the "0" / "***.*" placeholder text ---
# These use the same shared-string placeholders already used elsewhere in the
# sheet (e.g. A14/C14 style). Writing the text first (apostrophe-prefixed so a
# value that looks numeric, like "0", is not auto-converted back to a number),
# then pasting the number-format/style from a cell that already carries the
# placeholder style (style index 14) keeps formatting consistent with the
# rest of the column instead of leaving a stray "Text" number format behind.
$ws.Range("C15").Value = "'0"
$ws.Range("D15").Value = "'0"
$ws.Range("E15").Value = "***.*"
$ws.Range("D26").Value = "'0"
$ws.Range("E26").Value = "***.*"

$ws.Range("A15").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").PasteSpecial(-4122)

$ws.Range("A26").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("E26").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- New crime-data figures for rows 15-21, 24-27, 28-29 ---
$ws.Range("M15").Value = 10
$ws.Range("N15").Value = -45
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -20
$ws.Range("F16").Value = 17
$ws.Range("G16").Value = 20
$ws.Range("H16").Value = -15
$ws.Range("I16").Value = 114
$ws.Range("J16").Value = 138
$ws.Range("K16").Value = -17.391304347826
$ws.Range("L16").Value = 54.054054054054
$ws.Range("M16").Value = -37.016574585635
$ws.Range("N16").Value = -84.820239680426
$ws.Range("C17").Value = 10
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = 42.857142857142
$ws.Range("G17").Value = 29
$ws.Range("H17").Value = 27.586206896551
$ws.Range("I17").Value = 255
$ws.Range("J17").Value = 257
$ws.Range("K17").Value = -0.778210116731
$ws.Range("L17").Value = 45.714285714285
$ws.Range("M17").Value = 65.584415584415
$ws.Range("N17").Value = 12.831858407079
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = -80
$ws.Range("F18").Value = 8
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = -33.333333333333
$ws.Range("I18").Value = 96
$ws.Range("J18").Value = 93
$ws.Range("K18").Value = 3.225806451612
$ws.Range("L18").Value = 29.729729729729
$ws.Range("M18").Value = -53.623188405797
$ws.Range("N18").Value = -89.666307857911
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = 42.857142857142
$ws.Range("F19").Value = 24
$ws.Range("G19").Value = 28
$ws.Range("H19").Value = -14.285714285714
$ws.Range("I19").Value = 204
$ws.Range("J19").Value = 236
$ws.Range("K19").Value = -13.559322033898
$ws.Range("L19").Value = 31.612903225806
$ws.Range("M19").Value = -5.116279069767
$ws.Range("N19").Value = -42.696629213483
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = 16.666666666666
$ws.Range("F20").Value = 29
$ws.Range("G20").Value = 21
$ws.Range("H20").Value = 38.095238095238
$ws.Range("I20").Value = 188
$ws.Range("J20").Value = 165
$ws.Range("K20").Value = 13.939393939393
$ws.Range("L20").Value = 24.503311258278
$ws.Range("M20").Value = 8.670520231213
$ws.Range("N20").Value = -91.633288829550
$ws.Range("C21").Value = 32
$ws.Range("D21").Value = 30
$ws.Range("E21").Value = 6.666666666666
$ws.Range("F21").Value = 116
$ws.Range("G21").Value = 111
$ws.Range("H21").Value = 4.504504504504
$ws.Range("I21").Value = 872
$ws.Range("J21").Value = 905
$ws.Range("K21").Value = -3.646408839779
$ws.Range("L21").Value = 34.360554699537
$ws.Range("M21").Value = -7.529162248144
$ws.Range("N21").Value = -80.763291418486
$ws.Range("C24").Value = 19
$ws.Range("D24").Value = 49
$ws.Range("E24").Value = -61.224489795918
$ws.Range("F24").Value = 103
$ws.Range("G24").Value = 157
$ws.Range("H24").Value = -34.394904458598
$ws.Range("I24").Value = 971
$ws.Range("J24").Value = 1034
$ws.Range("K24").Value = -6.092843326885
$ws.Range("L24").Value = 56.109324758842
$ws.Range("M24").Value = 91.518737672583
$ws.Range("C25").Value = 18
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = 100
$ws.Range("F25").Value = 49
$ws.Range("G25").Value = 38
$ws.Range("H25").Value = 28.947368421052
$ws.Range("I25").Value = 361
$ws.Range("J25").Value = 354
$ws.Range("K25").Value = 1.977401129943
$ws.Range("L25").Value = 9.726443768996
$ws.Range("M25").Value = -13.012048192771
$ws.Range("F26").Value = 4
$ws.Range("H26").Value = 300
$ws.Range("I26").Value = 24
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = -7.692307692307
$ws.Range("D27").Value = 1
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = -25
$ws.Range("J27").Value = 35
$ws.Range("K27").Value = 20
$ws.Range("L27").Value = -6.666666666666
$ws.Range("N28").Value = -68.181818181818
$ws.Range("N29").Value = -72.222222222222
